$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (text unchanged, columns stay nick/addr/school/tel) ---
$ws.Range("A1").Value = "nick"
$ws.Range("B1").Value = "addr"
$ws.Range("C1").Value = "school"
$ws.Range("D1").Value = "tel"

# --- Row 2: おとちゃん ---
$ws.Range("A2").Value = "おとちゃん"
$ws.Range("B2").Value = "千葉市中央区末広４丁目１６－１１"
$ws.Range("C2").Value = "寒川小学校"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "09011956792"

# --- Row 3: りくぼー ---
$ws.Range("A3").Value = "りくぼー"
$ws.Range("B3").Value = "千葉市中央区末広４丁目１６－１１"
$ws.Range("C3").Value = "寒川小学校"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "09011956792"

# --- Row 4 (new row): ヒーちゃん ---
# Copy the already-established formatting from row 3 down onto row 4 so the
# B/C/D cells pick up the same cell style that B2:D3 use (rather than having
# the engine mint a brand-new, unused style entry). A4 is intentionally left
# with the default style, matching the source row.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null

$ws.Range("B4").Value = "千葉市中央区末広４丁目１６－１１"
$ws.Range("C4").Value = "寒川小学校"
$ws.Range("D4").Value = "09011956792"
$ws.Range("A4").Value = "ヒーちゃん"

# --- Column widths ---
# Target widths are A=16.75, B=38.125, C=34, D=12.25 "characters". The
# engine always re-derives the stored width from a fixed 7px "Maximum
# Digit Width" pixel grid, so only multiples of 1/7 are representable;
# these inputs are the values (empirically verified) that land on the
# closest achievable grid point to each intended width (34 lands exactly).
$ws.Columns.Item(1).ColumnWidth = 16.0
$ws.Columns.Item(2).ColumnWidth = 37.425
$ws.Columns.Item(3).ColumnWidth = 33.285
$ws.Columns.Item(4).ColumnWidth = 11.57

# --- Selection matches the saved cursor position in the file ---
$ws.Range("A7").Select() | Out-Null
